# Add a new student (Aawaj Joshi, s521315) to the "Student" sheet and the
# corresponding five StudentTerm rows to the "StudentTerm" sheet, then leave
# the UI selection/active-sheet state the way the edited workbook shows it
# (Student sheet active, cell F4 selected there; StudentTerm scrolled down
# with B17 selected).

$wb = $excel.ActiveWorkbook

# ---- Student sheet: new row for Aawaj Joshi (s521315) ----
$wsStudent = $wb.Worksheets.Item("Student")
$wsStudent.Range("A4").Value = 521315
$wsStudent.Range("B4").Value = "Aawaj"
$wsStudent.Range("C4").Value = "Joshi"
$wsStudent.Range("D4").Value = "s521315"
$wsStudent.Range("E4").Value = 480684

# ---- StudentTerm sheet: five new term rows for the new student ----
$wsStudentTerm = $wb.Worksheets.Item("StudentTerm")

$wsStudentTerm.Range("A13").Value = 12
$wsStudentTerm.Range("B13").Value = 521315
$wsStudentTerm.Range("C13").Value = 12
$wsStudentTerm.Range("D13").Value = "Spring 2018"

$wsStudentTerm.Range("A14").Value = 13
$wsStudentTerm.Range("B14").Value = 521315
$wsStudentTerm.Range("C14").Value = 13
$wsStudentTerm.Range("D14").Value = "Fall 2018"

$wsStudentTerm.Range("A15").Value = 14
$wsStudentTerm.Range("B15").Value = 521315
$wsStudentTerm.Range("C15").Value = 14
$wsStudentTerm.Range("D15").Value = "Spring 2019"

$wsStudentTerm.Range("A16").Value = 15
$wsStudentTerm.Range("B16").Value = 521315
$wsStudentTerm.Range("C16").Value = 15
$wsStudentTerm.Range("D16").Value = "Summer 2019"

$wsStudentTerm.Range("A17").Value = 16
$wsStudentTerm.Range("B17").Value = 521315
$wsStudentTerm.Range("C17").Value = 16
$wsStudentTerm.Range("D17").Value = "Fall 2019"

# ---- View state ----
# StudentTerm: scrolled so row 7 is at the top, B17 selected (becomes the
# active sheet momentarily while we set its selection).
[void]$wsStudentTerm.Select()
$excel.ActiveWindow.ScrollRow = 7
[void]$wsStudentTerm.Range("B17").Select()

# Student: ends up the active/selected tab with F4 selected.
[void]$wsStudent.Select()
[void]$wsStudent.Range("F4").Select()
